$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 42608.892118055555
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"

$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 56
$ws.Range("D7").Value = 37
$ws.Range("E7").Value = 47
$ws.Range("F7").Value = 52
$ws.Range("G7").Value = 20519
$ws.Range("H7").Value = 20679
$ws.Range("I7").Value = 3523
$ws.Range("J7").Value = 378
$ws.Range("K7").Value = 251
$ws.Range("L7").Value = 17
$ws.Range("M7").Value = 19
$ws.Range("N7").Value = "Noun"
